$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 2023-09-01 (45170) to 2023-09-05 (45174) for rows 2-7
$ws.Range("C2:C7").Value = 45174
